$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall Results")

# Activate this sheet (moves tabSelected from the previously active sheet)
$ws.Activate()

# Set new header cell values (order controls shared-string table order)
$ws.Range("A1").Value = "big column"
$ws.Range("C1").Value = "normal col"
$ws.Range("B1").Value = "tiny"

# Wrap text on the "tiny" header cell
$ws.Range("B1").WrapText = $true

# Adjust row height for the header row
$ws.Rows.Item(1).RowHeight = 48.75

# Adjust column widths (closest achievable values given the host's
# pixel-width rounding; targets are 32.28515625 and 4.42578125 chars)
$ws.Columns.Item(1).ColumnWidth = 31.5
$ws.Columns.Item(2).ColumnWidth = 3.67

# Update selection on the sheet
$ws.Range("D4").Select()
